$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of J2:J11
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14-17: summary labels and formulas
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Apply bold 12pt font with vertical-center alignment to B14:B17.
# Format B14 directly, then copy/paste the format onto B15:B17 so the
# whole block resolves to a single shared style (avoids minting one new
# style per cell).
$first = $ws.Range("B14")
$first.Font.Bold = $true
$first.Font.Size = 12
$first.VerticalAlignment = -4108

$first.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights for 14-17
$ws.Range("A14:A17").EntireRow.RowHeight = 15.6

# Selection
$ws.Range("A14:B17").Select()
